$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -8.4
$ws.Range("E12").Value = "71.6/140"
